$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the species/spicule label in A2
$ws.Range("A2").Value = "oxyaster euaster (bendy spines)"

# Clear row 3, which previously held empty placeholder cells, to be truly empty
$ws.Range("A3:F3").ClearContents()

# Update selection to match the saved workbook state
$ws.Range("E10").Select() | Out-Null
